# Auto-generated Excel COM-interop script to apply the commit diff
$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value2 = '2025-06-13T15:45:04+00:00'
$wsMeta.Range("B15").Value2 = '4.0.1'

# --- Elements sheet updates ---
$wsElem = $wb.Worksheets.Item("Elements")
$wsElem.Range("AJ2").Value2 = 'dom-2:If the resource is contained in another resource, it SHALL NOT contain nested Resources {contained.contained.empty()}
dom-3:If the resource is contained in another resource, it SHALL be referred to from elsewhere in the resource or SHALL refer to the containing resource {contained.where(((''#''+id in (%resource.descendants().reference | %resource.descendants().as(canonical) | %resource.descendants().as(uri) | %resource.descendants().as(url))) or descendants().where(reference = ''#'').exists() or descendants().where(as(canonical) = ''#'').exists() or descendants().where(as(canonical) = ''#'').exists()).not()).trace(''unmatched'', id).empty()}dom-4:If a resource is contained in another resource, it SHALL NOT have a meta.versionId or a meta.lastUpdated {contained.meta.versionId.empty() and contained.meta.lastUpdated.empty()}dom-5:If a resource is contained in another resource, it SHALL NOT have a security label {contained.meta.security.empty()}dom-6:A resource should have narrative for robust management {text.`div`.exists()}'
$wsElem.Range("AL2").Value2 = ''
$wsElem.Range("Y6").Value2 = 'A human language.'
$wsElem.Range("AJ8").Value2 = ''
$wsElem.Range("O10").Value2 = 'Modifier extensions allow for extensions that *cannot* be safely ignored to be clearly distinguished from the vast majority of extensions which can be safely ignored.  This promotes interoperability by eliminating the need for implementers to prohibit the presence of extensions. For further information, see the [definition of modifier extensions](http://hl7.org/fhir/R4/extensibility.html#modifierExtension).'
$wsElem.Range("Z14").Value2 = 'http://hl7.org/fhir/ValueSet/administrative-gender|4.0.1'
$wsElem.Range("AJ20").Value2 = ''
$wsElem.Range("O23").Value2 = 'Modifier extensions allow for extensions that *cannot* be safely ignored to be clearly distinguished from the vast majority of extensions which can be safely ignored.  This promotes interoperability by eliminating the need for implementers to prohibit the presence of extensions. For further information, see the [definition of modifier extensions](http://hl7.org/fhir/R4/extensibility.html#modifierExtension).'
$wsElem.Range("Z25").Value2 = 'http://hl7.org/fhir/ValueSet/identity-assuranceLevel|4.0.1'

Write-Host "Done applying edits."
